$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数" / want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10
$ws1.Range("F3").Value = 1045
$ws1.Range("F4").Value = 519
$ws1.Range("F5").Value = 13735
$ws1.Range("F7").Value = 41
$ws1.Range("F8").Value = 1760
$ws1.Range("F11").Value = 88
$ws1.Range("F12").Value = 46
$ws1.Range("F15").Value = 13761
$ws1.Range("F16").Value = 352
$ws1.Range("F17").Value = 612
$ws1.Range("F18").Value = 9065
$ws1.Range("F20").Value = 8176
$ws1.Range("F21").Value = 262
$ws1.Range("F26").Value = 158
$ws1.Range("F27").Value = 10
$ws1.Range("F28").Value = 26
$ws1.Range("F30").Value = 14
$ws1.Range("F33").Value = 6
$ws1.Range("F35").Value = 206
$ws1.Range("F38").Value = 5038

# Sheet "全部类型" (sheet4) updates to column F ("想去人数" / want-to-go count)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10
$ws4.Range("F3").Value = 1045
$ws4.Range("F4").Value = 519
$ws4.Range("F5").Value = 13735
$ws4.Range("F7").Value = 41
$ws4.Range("F8").Value = 1760
$ws4.Range("F11").Value = 88
$ws4.Range("F12").Value = 46
$ws4.Range("F15").Value = 13762
$ws4.Range("F16").Value = 352
$ws4.Range("F17").Value = 612
$ws4.Range("F18").Value = 9065
$ws4.Range("F20").Value = 8176
$ws4.Range("F21").Value = 262
$ws4.Range("F26").Value = 158
$ws4.Range("F27").Value = 10
$ws4.Range("F28").Value = 26
$ws4.Range("F30").Value = 14
$ws4.Range("F35").Value = 6
$ws4.Range("F37").Value = 206
$ws4.Range("F40").Value = 5038
